# ozon fixes 30.10.2025 part 2
#
# The first two orders (2056204 / 2059046) are removed from the sheet, so
# row 2 and row 3 are deleted outright and everything below shifts up two
# rows (the former row 4, order 2083871, becomes the new row 2; the former
# blank row 5 becomes row 3; etc. - and the trailing blank rows 21/22 go
# away, shrinking the used range from E22 down to E20).
#
# The order that lands on the new row 2 (2083871) also gets its to_cred /
# from_merch figures corrected down to 0, and the duplicate-value
# highlighting that used to flag the "order" column is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two fulfilled/removed orders - everything below (data + blank,
# pre-formatted rows) shifts up by two rows, matching dimension A1:E20.
$ws.Rows("2:3").Delete()

# The order that is now on row 2 had its credit/merchant adjustments
# corrected to zero.
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# The "highlight duplicate order numbers" conditional formatting on column A
# is no longer needed.
$ws.Range("A2").FormatConditions.Delete()

# Leave the selection where the user last clicked.
$ws.Range("F11").Select()
